# Generate Report for Archive
#
# The localization run for 3d910c7f-...md and bd85d307-...md moved from
# "Ready for handoff" into "In Translation", while 224bc1f8-...md's row
# (previously reported first) now reports last. This reshuffles rows
# 3-5 on every sheet (Overview / zh-cn / de-de) and refreshes the three
# hyperlink columns that reference the e2e file names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md"
$ws.Range("B3").Value = "e2e\3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md"
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"
$ws.Range("G3").Value = "2016-08-28 02:40:34"

$ws.Range("A4").Value = "bd85d307-11d0-4207-8838-4c6012e1889a.md"
$ws.Range("B4").Value = "e2e\bd85d307-11d0-4207-8838-4c6012e1889a.md"
$ws.Range("C4").Value = ".md"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"
$ws.Range("G4").Value = "2016-08-28 02:40:34"

$ws.Range("A5").Value = "224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md"
$ws.Range("B5").Value = "e2e\224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md"
$ws.Range("C5").Value = ".md"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-28 02:39:32"

# Rebuild the hyperlinks so each row's display text matches its (now
# shuffled) file name, while the underlying target URLs stay pinned to
# the same row position they were on before (rId2..rId5 keep pointing
# at the same commit URLs as in the original workbook).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md", "", "", "e2e\3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5970467c1098d4409e2b37952e381f3b98f6e23/e2e/224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md", "", "", "e2e\3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d822c47bc72984f86a3d1bec2df2ac0dbaedcb93/e2e/3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md", "", "", "e2e\bd85d307-11d0-4207-8838-4c6012e1889a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d822c47bc72984f86a3d1bec2df2ac0dbaedcb93/e2e/bd85d307-11d0-4207-8838-4c6012e1889a.md", "", "", "e2e\224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("G3").Value = "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.d62eaa93f590c0d247a28b32abb24731821c9c8b.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-28 02:40:29"

$ws.Range("A4").Value = "bd85d307-11d0-4207-8838-4c6012e1889a.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "bd85d307-11d0-4207-8838-4c6012e1889a.ded0bf512234fe44a8e2b2ec0e81482c79ea91bb.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-28 02:40:29"

$ws.Range("A5").Value = "224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("G5").Value = "224bc1f8-8fbf-4d73-b0a5-ce37facc811b.d6b9d11cad8957a16827a73521f6084c54b8e61e.zh-cn.xlf"
$ws.Range("H5").Value = "2016-08-28 02:39:28"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md", "", "", "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/debd59d6ea1a13b893df256aea5f697faf7b79a2/e2e/3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md", "", "", "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5970467c1098d4409e2b37952e381f3b98f6e23/e2e/224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md", "", "", "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d822c47bc72984f86a3d1bec2df2ac0dbaedcb93/e2e/3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md", "", "", "bd85d307-11d0-4207-8838-4c6012e1889a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d822c47bc72984f86a3d1bec2df2ac0dbaedcb93/e2e/bd85d307-11d0-4207-8838-4c6012e1889a.md", "", "", "224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md"
$ws.Range("C3").Value = "In Translation"
$ws.Range("G3").Value = "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.d62eaa93f590c0d247a28b32abb24731821c9c8b.de-de.xlf"
$ws.Range("H3").Value = "2016-08-28 02:40:34"

$ws.Range("A4").Value = "bd85d307-11d0-4207-8838-4c6012e1889a.md"
$ws.Range("C4").Value = "In Translation"
$ws.Range("G4").Value = "bd85d307-11d0-4207-8838-4c6012e1889a.ded0bf512234fe44a8e2b2ec0e81482c79ea91bb.de-de.xlf"
$ws.Range("H4").Value = "2016-08-28 02:40:34"

$ws.Range("A5").Value = "224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("G5").Value = "224bc1f8-8fbf-4d73-b0a5-ce37facc811b.d6b9d11cad8957a16827a73521f6084c54b8e61e.de-de.xlf"
$ws.Range("H5").Value = "2016-08-28 02:39:32"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4401b0514e73b49f6519f4fa733bcdb83bd96d7d/e2e/3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md", "", "", "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c6f0bea01bedace84758765734320bad68174b8d/e2e/3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md", "", "", "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b5970467c1098d4409e2b37952e381f3b98f6e23/e2e/224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md", "", "", "3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d822c47bc72984f86a3d1bec2df2ac0dbaedcb93/e2e/3d910c7f-3e0d-4ff9-b600-e52dbdb56ac6.md", "", "", "bd85d307-11d0-4207-8838-4c6012e1889a.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d822c47bc72984f86a3d1bec2df2ac0dbaedcb93/e2e/bd85d307-11d0-4207-8838-4c6012e1889a.md", "", "", "224bc1f8-8fbf-4d73-b0a5-ce37facc811b.md") | Out-Null
